# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" sheets to reflect the newly generated data snapshot.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8985
    $ws.Range("F4").Value = 448
}
